$wb = $excel.ActiveWorkbook

# Worksheet "展览" (sheet1): update "想去人数" (column F) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 2002
$wsExhibit.Range("F5").Value = 316
$wsExhibit.Range("F6").Value = 64
$wsExhibit.Range("F8").Value = 2046
$wsExhibit.Range("F9").Value = 10385
$wsExhibit.Range("F14").Value = 396
$wsExhibit.Range("F15").Value = 7236
$wsExhibit.Range("F18").Value = 140
$wsExhibit.Range("F20").Value = 269

# Worksheet "全部类型" (sheet4): update "想去人数" (column F) values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 2002
$wsAll.Range("F5").Value = 316
$wsAll.Range("F6").Value = 64
$wsAll.Range("F9").Value = 2046
$wsAll.Range("F12").Value = 10385
$wsAll.Range("F17").Value = 396
$wsAll.Range("F18").Value = 7236
$wsAll.Range("F21").Value = 140
$wsAll.Range("F23").Value = 269
